$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert two new sheets: "AppName" and "ModuleName", right after
#    "AddOpportunity" and before "Users".
# ---------------------------------------------------------------
$appNameSheet = $wb.Worksheets.Add()
$appNameSheet.Name = "AppName"
$appNameSheet.Move($null, $wb.Worksheets.Item("AddOpportunity"))

$moduleNameSheet = $wb.Worksheets.Add()
$moduleNameSheet.Name = "ModuleName"
$moduleNameSheet.Move($null, $wb.Worksheets.Item("AppName"))

# ---------------------------------------------------------------
# 2. Populate "AppName" sheet
# ---------------------------------------------------------------
$appName = $wb.Worksheets.Item("AppName")
$appName.Range("A1").Value = "App Name"
$appName.Range("A1").Font.Bold = $true
$appName.Range("A2").Value = "HL Banker"
$appName.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------
# 3. Populate "ModuleName" sheet
# ---------------------------------------------------------------
$moduleName = $wb.Worksheets.Item("ModuleName")
$moduleName.Range("A1").Value = "Module Name"
$moduleName.Range("A1").Font.Bold = $true
$moduleName.Range("A3").Value = "Cases"
$moduleName.Range("A2").Value = "Opportunities"
$moduleName.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------
# 4. AddContact sheet: selection moves from D14 to A2
# ---------------------------------------------------------------
$addContact = $wb.Worksheets.Item("AddContact")
$addContact.Range("A2").Select()

# ---------------------------------------------------------------
# 5. OppDealTeamMembers sheet: clear the leftover font-style on the
#    three "Specialty" section-header rows & move the selection.
# ---------------------------------------------------------------
$oppDealTeam = $wb.Worksheets.Item("OppDealTeamMembers")
$oppDealTeam.Range("A13").ClearFormats()
$oppDealTeam.Range("A19").ClearFormats()
$oppDealTeam.Range("A28").ClearFormats()
$oppDealTeam.Range("A17").Select()

# ---------------------------------------------------------------
# 6. EngDealTeamMembers sheet: clear the leftover font-style on the
#    last row.
# ---------------------------------------------------------------
$engDealTeam = $wb.Worksheets.Item("EngDealTeamMembers")
$engDealTeam.Range("A2").ClearFormats()

# ---------------------------------------------------------------
# 7. AddOpportunity sheet becomes the active / selected sheet again,
#    with a new selection.
# ---------------------------------------------------------------
$addOpportunity = $wb.Worksheets.Item("AddOpportunity")
$addOpportunity.Activate()
$addOpportunity.Range("K2").Select()
